# Add a "Croatia" market test-data sheet, cloned from the existing
# "Turkey" sheet (same layout/styles), populated with Croatia-specific
# values, and placed as the new last tab (made active).

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Mimic the "select everything" state Turkey ends up with once it is no
# longer the focused tab.
$turkey.Activate()
$null = $turkey.Cells.Select()

# Clone Turkey's sheet (keeps column widths, styles, merged cells, etc.)
# and drop it immediately after Turkey -> becomes the last tab.
$turkey.Copy($null, $turkey)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Market-specific values.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2473"

# Leave the cursor on B4 and make Croatia the active sheet/tab.
$null = $croatia.Range("B4").Select()
$croatia.Activate()
